$wb = $excel.ActiveWorkbook

$weeks = 1..15

foreach ($week in $weeks) {
    $ws = $wb.Worksheets.Item("GV_Tuan_$week")

    # ---- Row 8: "C2 (15:00-17:00)" slot becomes "T1 (17:30-19:30)" ----
    # The class that used to sit in C8 (Monday / "Thu 2") moves to D8 (Tuesday / "Thu 3"),
    # keeping its highlighted style; C8 goes back to being a plain empty slot.
    $ws.Range("A8").Value = "T1`n(17:30-19:30)"

    $ws.Range("C8").Copy()
    $ws.Range("D8").PasteSpecial(-4122)
    $ws.Range("D8").Value = "Lớp: CL10`nMôn: Ngữ pháp tiếng Anh`nPhòng: R105`n(Lý thuyết)"

    $ws.Range("F8").Copy()
    $ws.Range("C8").PasteSpecial(-4122)
    $ws.Range("C8").Value = ""

    # ---- Row 9: "T1 (17:30-19:30)" slot becomes "T2 (19:30-21:30)" ----
    # The class that used to sit in E9 (Wednesday / "Thu 4") moves to C9 (Monday / "Thu 2"),
    # keeping its highlighted style, and its room changes from R102 to R103;
    # E9 goes back to being a plain empty slot.
    $ws.Range("A9").Value = "T2`n(19:30-21:30)"

    $ws.Range("E9").Copy()
    $ws.Range("C9").PasteSpecial(-4122)
    $ws.Range("C9").Value = "Lớp: CL05`nMôn: Ngữ pháp tiếng Anh`nPhòng: R103`n(Lý thuyết)"

    $ws.Range("F9").Copy()
    $ws.Range("E9").PasteSpecial(-4122)
    $ws.Range("E9").Value = ""
}

$excel.CutCopyMode = 0
